# "Fruta / hortaliza, semanal" - weekly update.
# A new weekly record is inserted as row 68 (pushing the existing rows 68-106
# down to 69-107), and the sheet's used-range dimension grows from R106 to R107.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 68; existing rows 68..106 shift to 69..107.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new weekly price record.
$ws.Range("A68").Value = 3
$ws.Range("B68").Value = "Femacal de La Calera"
$ws.Range("C68").Value = "Coquimbo"
$ws.Range("D68").Value = 44529
$ws.Range("E68").Value = 5
$ws.Range("F68").Value = 100112030
$ws.Range("G68").Value = "Poroto granado"
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 43
$ws.Range("K68").Value = 34000
$ws.Range("L68").Value = 35000
$ws.Range("M68").Value = 34465
$ws.Range("N68").Value = "$/malla 25 kilos"
$ws.Range("O68").Value = "Provincia de Limarí"
$ws.Range("P68").Value = 1379
$ws.Range("Q68").Value = 25
$ws.Range("R68").Value = "Hortaliza"
